$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B81 to be a numeric value instead of an inline string "4"
$ws.Range("B81").Value = 4

# Add new row 82 with the additional annotation data
$ws.Range("A82").Value = "Ying Tang"

# politeness_score (B82) stays a text value "3" (matches source data style of most rows,
# but keep it as text like the original diff expects)
$ws.Range("B82").NumberFormat = "@"
$ws.Range("B82").Value = "3"
$ws.Range("B82").Style = "Normal"

$ws.Range("C82").Value = "无"
$ws.Range("D82").Value = "DFT"
$ws.Range("E82").Value = "WRI"
$ws.Range("F82").Value = "a443a511-4201-4c2e-8ed3-4de4ebdb5fb0"
$ws.Range("G82").Value = "BkfEzz-0-_annotated.xlsx"
$ws.Range("H82").Value = "There are several important concepts, such as reward distribution, credit assignment, which are used (from the very beginning of the paper) without explanation until the final part of the paper."
